# Add a new "Source" column (AM) to the "List of Quals" sheet:
#   - AM3 header = "Source" (same header style as the adjoining AL3 cell)
#   - AM4 value  = "Gatsby"
# and leave the newly added cell selected, matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell, column AM (39), row 3 - reuse the formatting of the
# last existing header cell (AL3) so the rotated/shaded style carries over.
$ws.Range("AM3").Value = "Source"
$ws.Range("AL3").Copy()
$ws.Range("AM3").PasteSpecial(-4122)  # xlPasteFormats

# New data cell, column AM (39), row 4.
$ws.Range("AM4").Value = "Gatsby"

# Match the final selection left by the author's edit.
$ws.Activate() | Out-Null
$ws.Range("AM4").Select() | Out-Null
